$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 357.72726
$ws.Range("I58").Value = 357.72726
$ws.Range("K58").Value = 1073.18178
$ws.Range("M58").Value = -923.1817799999999
$ws.Range("H64").Value = 4500
$ws.Range("H67").Value = 4500
$ws.Range("H69").Value = 7999.75
$ws.Range("I69").Value = 7333
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 21999
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -21125
$ws.Range("N69").Value = -31748
$ws.Range("H72").Value = 7999.75
$ws.Range("I72").Value = 7333
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 65997
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -61629
$ws.Range("N72").Value = -98736
$ws.Range("H87").Value = 12607677
$ws.Range("J87").Value = 12607677
$ws.Range("L87").Value = 12607677
$ws.Range("N87").Value = -12610173
$ws.Range("H90").Value = 12607677
$ws.Range("J90").Value = 12607677
$ws.Range("L90").Value = 37823031
$ws.Range("N90").Value = -37835511
$ws.Range("H112").Value = 3586.0815
$ws.Range("J112").Value = 3602.8262
$ws.Range("L112").Value = 10808.4786
$ws.Range("N112").Value = -13024.4786
$ws.Range("H132").Value = 2059.9697
$ws.Range("I132").Value = 2060.8923
$ws.Range("K132").Value = 6182.6769
$ws.Range("M132").Value = -3652.6769
$ws.Range("H137").Value = 3016
$ws.Range("J137").Value = 2073.3333
$ws.Range("L137").Value = 6219.999899999999
$ws.Range("N137").Value = -11319.9999
$ws.Range("H138").Value = 152440.81
$ws.Range("I138").Value = 5958.5454
$ws.Range("J138").Value = 180221.94
$ws.Range("K138").Value = 17875.6362
$ws.Range("L138").Value = 540665.8200000001
$ws.Range("M138").Value = -12735.6362
$ws.Range("N138").Value = -550945.8200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2537.074
$ws.Range("I45").Value = 2004.4706
$ws.Range("J45").Value = 3442.5
$ws.Range("K45").Value = 2004.4706
$ws.Range("L45").Value = 3442.5
$ws.Range("M45").Value = -1627.4706
$ws.Range("N45").Value = -4196.5
$ws.Range("H132").Value = 2082391.1
$ws.Range("I132").Value = 2655132.5
$ws.Range("K132").Value = 7965397.5
$ws.Range("M132").Value = -7962867.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 60000
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61622
$ws.Range("H71").Value = 60000
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -188112
$ws.Range("H86").Value = 5589222
$ws.Range("I86").Value = 5589222
$ws.Range("K86").Value = 5589222
$ws.Range("M86").Value = -5588099
$ws.Range("H87").Value = 147450
$ws.Range("J87").Value = 147450
$ws.Range("L87").Value = 147450
$ws.Range("N87").Value = -149946
$ws.Range("H89").Value = 5589222
$ws.Range("I89").Value = 5589222
$ws.Range("K89").Value = 27946110
$ws.Range("M89").Value = -27940494
$ws.Range("H90").Value = 147450
$ws.Range("J90").Value = 147450
$ws.Range("L90").Value = 442350
$ws.Range("N90").Value = -454830
$ws.Range("H99").Value = 3436.4583
$ws.Range("I99").Value = 2640.1177
$ws.Range("J99").Value = 5370.4287
$ws.Range("K99").Value = 2640.1177
$ws.Range("L99").Value = 5370.4287
$ws.Range("M99").Value = -1142.1177
$ws.Range("N99").Value = -8366.4287
$ws.Range("H134").Value = 3974.9412
$ws.Range("J134").Value = 4739.5
$ws.Range("L134").Value = 14218.5
$ws.Range("N134").Value = -19288.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3930.7097
$ws.Range("J31").Value = 4322.9165
$ws.Range("L31").Value = 4322.9165
$ws.Range("N31").Value = -4912.9165
$ws.Range("H34").Value = 3930.7097
$ws.Range("J34").Value = 4322.9165
$ws.Range("L34").Value = 4322.9165
$ws.Range("N34").Value = -4726.9165
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 52500
$ws.Range("J60").Value = 65000
$ws.Range("L60").Value = 65000
$ws.Range("N60").Value = -66022
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 9499.8
$ws.Range("J62").Value = 9333.333
$ws.Range("L62").Value = 9333.333
$ws.Range("N62").Value = -10581.333
$ws.Range("H65").Value = 9499.8
$ws.Range("J65").Value = 9333.333
$ws.Range("L65").Value = 46666.665
$ws.Range("N65").Value = -52906.665
$ws.Range("H68").Value = 65295
$ws.Range("J68").Value = 65295
$ws.Range("L68").Value = 65295
$ws.Range("N68").Value = -66793
$ws.Range("H71").Value = 65295
$ws.Range("J71").Value = 65295
$ws.Range("L71").Value = 195885
$ws.Range("N71").Value = -203373
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H94").Value = 1609.1177
$ws.Range("I94").Value = 346.2857
$ws.Range("J94").Value = 2493.1
$ws.Range("K94").Value = 346.2857
$ws.Range("L94").Value = 2493.1
$ws.Range("M94").Value = 104.7143
$ws.Range("N94").Value = -3395.1
$ws.Range("H132").Value = 4150.108
$ws.Range("I132").Value = 3117.5
$ws.Range("J132").Value = 6590.8184
$ws.Range("K132").Value = 9352.5
$ws.Range("L132").Value = 19772.4552
$ws.Range("M132").Value = -6822.5
$ws.Range("N132").Value = -24832.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 206.25
$ws.Range("I2").Value = 166.66667
$ws.Range("J2").Value = 230
$ws.Range("K2").Value = 1000.00002
$ws.Range("L2").Value = 1380
$ws.Range("M2").Value = -887.0000200000001
$ws.Range("N2").Value = -1606
$ws.Range("H5").Value = 1457.5
$ws.Range("I5").Value = 636.25
$ws.Range("K5").Value = 1908.75
$ws.Range("M5").Value = -1796.75
$ws.Range("H28").Value = 2528.3333
$ws.Range("I28").Value = 2126
$ws.Range("K28").Value = 6378
$ws.Range("M28").Value = -6146
$ws.Range("H45").Value = 4044.3333
$ws.Range("J45").Value = 4044.3333
$ws.Range("L45").Value = 12132.9999
$ws.Range("N45").Value = -13196.9999
$ws.Range("H70").Value = 6270.3335
$ws.Range("J70").Value = 6270.3335
$ws.Range("L70").Value = 18811.0005
$ws.Range("N70").Value = -19441.0005
$ws.Range("H73").Value = 6270.3335
$ws.Range("J73").Value = 6270.3335
$ws.Range("L73").Value = 18811.0005
$ws.Range("N73").Value = -20995.0005
$ws.Range("H80").Value = 6588.5386
$ws.Range("J80").Value = 6643.625
$ws.Range("L80").Value = 19930.875
$ws.Range("N80").Value = -21802.875
$ws.Range("H83").Value = 6588.5386
$ws.Range("J83").Value = 6643.625
$ws.Range("L83").Value = 59792.625
$ws.Range("N83").Value = -69152.625
$ws.Range("H113").Value = 3499.1538
$ws.Range("J113").Value = 4398.778
$ws.Range("L113").Value = 13196.334
$ws.Range("N113").Value = -17536.334
$ws.Range("H131").Value = 1252622.6
$ws.Range("J131").Value = 2997.4285
$ws.Range("L131").Value = 8992.2855
$ws.Range("N131").Value = -19072.2855
$ws.Range("H132").Value = 3292
$ws.Range("J132").Value = 3446.5
$ws.Range("L132").Value = 31018.5
$ws.Range("N132").Value = -36078.5
$ws.Range("H135").Value = 1457.5
$ws.Range("I135").Value = 636.25
$ws.Range("K135").Value = 5726.25
$ws.Range("M135").Value = -3191.25
$ws.Range("H141").Value = 7519.933
$ws.Range("I141").Value = 5414.2144
$ws.Range("K141").Value = 16242.6432
$ws.Range("M141").Value = -11062.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24111.938
$ws.Range("J57").Value = 69992.5
$ws.Range("L57").Value = 69992.5
$ws.Range("N57").Value = -71632.5
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 4621.434
$ws.Range("I132").Value = 4941.0444
$ws.Range("J132").Value = 2823.625
$ws.Range("K132").Value = 14823.1332
$ws.Range("L132").Value = 8470.875
$ws.Range("M132").Value = -12293.1332
$ws.Range("N132").Value = -13530.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16381.923
$ws.Range("I7").Value = 18107.889
$ws.Range("K7").Value = 18107.889
$ws.Range("M7").Value = -17995.889
$ws.Range("H46").Value = 6240.067
$ws.Range("I46").Value = 3130.4348
$ws.Range("J46").Value = 16457.428
$ws.Range("K46").Value = 3130.4348
$ws.Range("L46").Value = 16457.428
$ws.Range("M46").Value = -2942.4348
$ws.Range("N46").Value = -16833.428
$ws.Range("H68").Value = 1900
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 1900
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H126").Value = 16381.923
$ws.Range("I126").Value = 18107.889
$ws.Range("K126").Value = 54323.667
$ws.Range("M126").Value = -51853.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 398888.5
$ws.Range("I62").Value = 398888.5
$ws.Range("K62").Value = 398888.5
$ws.Range("M62").Value = -398264.5
$ws.Range("H65").Value = 398888.5
$ws.Range("I65").Value = 398888.5
$ws.Range("K65").Value = 1994442.5
$ws.Range("M65").Value = -1991322.5
$ws.Range("H81").Value = 56639.25
$ws.Range("J81").Value = 10749.5
$ws.Range("L81").Value = 21499
$ws.Range("N81").Value = -23621
$ws.Range("H84").Value = 56639.25
$ws.Range("J84").Value = 10749.5
$ws.Range("L84").Value = 107495
$ws.Range("N84").Value = -118103
